# Add a second strategy chart sheet ("Chart2") that mirrors the first
# ("Chart1"): same grid of values/labels and its own (placeholder) chart
# drawing — exactly what Excel does when you duplicate a sheet via
# "Move or Copy... > Create a copy".

$wb = $excel.ActiveWorkbook

# The existing (only) sheet holding the first strategy chart's data.
$ws1 = $wb.Worksheets.Item(1)

# Duplicate it, placing the copy immediately after Chart1. Excel gives the
# duplicate a name like "Chart1 (2)" and makes it the active sheet.
$ws1.Copy($null, $ws1)

# Rename the freshly created copy to "Chart2".
$ws2 = $wb.ActiveSheet
$ws2.Name = "Chart2"
